$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = "| HR260K63379m`n"
$ws.Range("A44").Value = "W720 BOM`n"

$ws.Rows.Item(43).AutoFit()
$ws.Rows.Item(44).AutoFit()
